$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string "og_group_ref" -> "field_collection_field" wherever it appears (A17)
$ws.Range("A17").Value = "field_collection_field"

# Update the selected/active cell in the sheet view from D18 to B17
$ws.Range("B17").Select()

# Update cell B17's value from 144795 to 1566
$ws.Range("B17").Value = 1566
